$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13: add model headers in L13/M13/N13
$ws.Range("L13").Value = "WMM15"
$ws.Range("M13").Value = "EMM15"
$ws.Range("N13").Value = "IGRF12"

# Row 14: "Model" / "WMM" labels, plus Bx model comparison row
$ws.Range("B14").Value = "Model"
$ws.Range("C14").Value = "WMM"
$ws.Range("K14").Value = "Bx"
$ws.Range("L14").Value = -34.8709
$ws.Range("M14").Value = -34.92592
$ws.Range("N14").Value = -34.911878

# Row 15: By model comparison row
$ws.Range("K15").Value = "By"
$ws.Range("L15").Value = 166.089
$ws.Range("M15").Value = 166.181764
$ws.Range("N15").Value = 166.096031

# Row 16: Bz model comparison row
$ws.Range("K16").Value = "Bz"
$ws.Range("L16").Value = 389.254
$ws.Range("M16").Value = 389.316472
$ws.Range("N16").Value = 389.318141

# Row 17: H model comparison row
$ws.Range("K17").Value = "H"
$ws.Range("L17").Value = 169.71
$ws.Range("M17").Value = 169.812245
$ws.Range("N17").Value = 169.725457

# Row 18: F model comparison row
$ws.Range("K18").Value = "F"
$ws.Range("L18").Value = 424.641
$ws.Range("M18").Value = 424.739348
$ws.Range("N18").Value = 424.706188

# Row 19: D model comparison row
$ws.Range("K19").Value = "D"
$ws.Range("L19").Value = -11.8572
$ws.Range("M19").Value = -11.868942
$ws.Range("N19").Value = -11.870255

# Row 20: I model comparison row
$ws.Range("K20").Value = "I"
$ws.Range("L20").Value = 66.4434
$ws.Range("M20").Value = 66.434109
$ws.Range("N20").Value = 66.444931

# Update selection to match the recorded final state
$ws.Range("G21").Select()
